# Refresh the cryptos price/volume snapshot (D/E columns) to the latest scrape.
# A leading apostrophe forces Excel to keep numeric-looking price strings as text
# (matching the source data, which stores prices/volumes as plain text, not numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.999.23"
$ws.Range("E2").Value = "  -0.43%  "

$ws.Range("D3").Value = "1.639.39"
$ws.Range("E3").Value = "  -1.05%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'214.93"
$ws.Range("E5").Value = "  -1.47%  "

$ws.Range("D6").Value = "'0.5051"
$ws.Range("E6").Value = "  -2.11%  "

$ws.Range("D7").Value = "'1.006"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "'0.06453"
$ws.Range("E8").Value = "  +0.36%  "

$ws.Range("D9").Value = "'0.2573"
$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("D10").Value = "'19.49"
$ws.Range("E10").Value = "  -1.91%  "

$ws.Range("D11").Value = "'0.07707"
$ws.Range("E11").Value = "  -0.62%  "

$ws.Range("D12").Value = "1.646.60"
$ws.Range("E12").Value = "  -1.32%  "

$ws.Range("D14").Value = "1.865.38"
$ws.Range("E14").Value = "  -0.99%  "

$ws.Range("D15").Value = "'0.5447"
$ws.Range("E15").Value = "  -1.61%  "

$ws.Range("D16").Value = "0.0₅7936"
$ws.Range("E16").Value = "  -1.31%  "

$ws.Range("D17").Value = "'63.37"
$ws.Range("E17").Value = "  -1.46%  "

$ws.Range("D18").Value = "25.981.26"
$ws.Range("E18").Value = "  -0.68%  "

$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").Value = "'203.41"
$ws.Range("E20").Value = "  -3.37%  "

$ws.Range("E21").Value = "  -2.39%  "

$ws.Range("D22").Value = "'9.999"
$ws.Range("E22").Value = "  -0.53%  "

$ws.Range("D23").Value = "'5.977"
$ws.Range("E23").Value = "  +1.03%  "

$ws.Range("D24").Value = "'1.007"
$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("D25").Value = "'1.959"
$ws.Range("E25").Value = "  +11.33%  "

$ws.Range("D26").Value = "'141.64"
$ws.Range("E26").Value = "  -1.54%  "

$ws.Range("E27").Value = "  -0.66%  "

$ws.Range("D28").Value = "'15.68"
$ws.Range("E28").Value = "  -0.30%  "

$ws.Range("D29").Value = "'6.720"
$ws.Range("E29").Value = "  -3.64%  "

$ws.Range("D30").Value = "'0.05043"

$ws.Range("D31").Value = "'1.239"
$ws.Range("E31").Value = "  -1.07%  "

$ws.Range("E32").Value = "  -3.05%  "

$ws.Range("D33").Value = "'3.193"
$ws.Range("E33").Value = "  -1.32%  "

$ws.Range("E34").Value = "  -2.10%  "

$ws.Range("E35").Value = "  -1.05%  "

$ws.Range("E36").Value = "  -4.36%  "

$ws.Range("D37").Value = "'0.8915"
$ws.Range("E37").Value = "  -3.41%  "

$ws.Range("D38").Value = "'0.5628"
$ws.Range("E38").Value = "  -0.93%  "

$ws.Range("D39").Value = "1.149.66"
$ws.Range("E39").Value = "  -1.26%  "

$ws.Range("E40").Value = "  -1.22%  "

$ws.Range("D41").Value = "'2.562"
$ws.Range("E41").Value = "  -0.15%  "

$ws.Range("D42").Value = "'1.006"
$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("E43").Value = "  +0.29%  "

$ws.Range("D44").Value = "'0.8080"
$ws.Range("E44").Value = "  -3.53%  "

$ws.Range("D45").Value = "'99.59"
$ws.Range("E45").Value = "  -0.30%  "

$ws.Range("D46").Value = "1.777.68"
$ws.Range("E46").Value = "  -0.92%  "

$ws.Range("E47").Value = "  +2.79%  "

$ws.Range("D48").Value = "'0.4524"
$ws.Range("E48").Value = "  +0.27%  "

$ws.Range("D49").Value = "'1.007"
$ws.Range("E49").Value = "  +0.05%  "

$ws.Range("D50").Value = "'54.97"
$ws.Range("E50").Value = "  -1.75%  "

$ws.Range("D51").Value = "'0.05032"
$ws.Range("E51").Value = "  -0.71%  "
